$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.004.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.417.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.47%  "

$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.002.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.85%  "

$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.420.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.75%  "

$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.017.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.20%  "

$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.87%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.555.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("E26").Value = "  -4.36%  "

$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -4.28%  "

$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("E33").Value = "  -2.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "168.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.449.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0778"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.776"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("E45").Value = "  -3.65%  "

$ws.Range("E46").Value = "  -5.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.544.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.56%  "

$ws.Range("E51").Value = "  -6.79%  "
